# "Generate Report for Handback"
#
# The handback-status report records, per locale, when each file was
# handed off for localization (Correspond Handoff Datetime) and when the
# corresponding handback file came back in sync (Correspond Handback
# DateTime). This run refreshes those timestamps for the
# 35e630c0-739c-4dd4-8d2e-e629d0444519 entry (row 2) on both the zh-cn and
# de-de report sheets. The Overview sheet's "Latest HO Xliff Generate
# Date" column (the max handoff datetime across locales for each file) is
# then re-generated for both files: it advances for
# 35e630c0-739c-4dd4-8d2e-e629d0444519.md (row 2, now driven by the new
# de-de handoff time) and is re-written unchanged for
# d6e383f1-9649-459a-a377-455aeec25a36.md (row 3).

$wb = $excel.ActiveWorkbook

# zh-cn sheet: update Correspond Handoff Datetime / Correspond Handback
# DateTime for the 35e630c0-739c-4dd4-8d2e-e629d0444519 row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-17 18:46:58"
$wsZhCn.Range("K2").Value = "2016-08-17 18:47:28"

# de-de sheet: update Correspond Handoff Datetime / Correspond Handback
# DateTime for the 35e630c0-739c-4dd4-8d2e-e629d0444519 row.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-17 18:47:08"
$wsDeDe.Range("K2").Value = "2016-08-17 18:47:35"

# Overview sheet: re-generate Latest HO Xliff Generate Date per file as
# the max Correspond Handoff Datetime across locales.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-17 18:47:08"
$wsOverview.Range("G3").Value = "2016-08-17 18:46:10"
